# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# The source feed re-ordered/re-matched a handful of fixture records for
# "Germany Regionalliga South West". The row *positions* (column A, the
# sequential index) stay put, but the bookmaker data attached to certain
# rows (id, teams, scores, odds, etc. - columns B through AD) needs to be
# reassigned to the row above (with wrap-around) for one block, and swapped
# between two rows for another block.
#
# Group 1: rows 257 and 260 simply trade places (B:AD).
# Group 2: rows 299-304 rotate down by one row, with row 304's data
#          wrapping around into row 299.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group 1: swap rows 257 and 260 (columns B:AD) ---
$row257 = $ws.Range("B257:AD257").Value2
$row260 = $ws.Range("B260:AD260").Value2

$ws.Range("B257:AD257").Value2 = $row260
$ws.Range("B260:AD260").Value2 = $row257

# --- Group 2: rotate rows 299-304 (columns B:AD) ---
# new(299) = old(304); new(300) = old(299); new(301) = old(300);
# new(302) = old(301); new(303) = old(302); new(304) = old(303)
$row299 = $ws.Range("B299:AD299").Value2
$row300 = $ws.Range("B300:AD300").Value2
$row301 = $ws.Range("B301:AD301").Value2
$row302 = $ws.Range("B302:AD302").Value2
$row303 = $ws.Range("B303:AD303").Value2
$row304 = $ws.Range("B304:AD304").Value2

$ws.Range("B299:AD299").Value2 = $row304
$ws.Range("B300:AD300").Value2 = $row299
$ws.Range("B301:AD301").Value2 = $row300
$ws.Range("B302:AD302").Value2 = $row301
$ws.Range("B303:AD303").Value2 = $row302
$ws.Range("B304:AD304").Value2 = $row303
